# UC04-Create Incident.docx edit:
# "Defines the incident and resources area on the map refer to "
#   -> "Defines the incident and resources location on the map refer to "
# with the new word "location" (and the following " on the map refer to ")
# ending up as their own separate runs (as produced by a real Word edit
# session), rather than merged back into the original run.

$d = $word.ActiveDocument

# --- Locate the word "area" inside the unique sentence and replace it ----
$anchor = $d.Content
$anchor.Find.Execute("resources area on the map refer to")
$sentenceStart = $anchor.Start

$oldWord = "area"
$newWord = "location"

$wordStart = $sentenceStart + "resources ".Length
$wordEnd = $wordStart + $oldWord.Length

$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Text = $newWord

# --- Force "location" to live in its own run (splits it away from the --
# --- text that follows/precedes it) by toggling a character property --
# --- on and back off again; this leaves no visible formatting change. --
$newWordRange = $d.Range($wordStart, $wordStart + $newWord.Length)
$newWordRange.Font.Bold = 1
$newWordRange.Font.Bold = 0

# --- The text replacement above also merges the unrelated, already- ----
# --- separate "Define Region on map " / "use-case" runs further along --
# --- in the same paragraph (they share identical formatting). Restore -
# --- that original run split the same way. --------------------------
$useCase = $d.Content
$useCase.Find.Execute("use-case")
$useCaseRange = $d.Range($useCase.Start, $useCase.End)
$useCaseRange.Font.Bold = 0
$useCaseRange.Font.Bold = 1
